$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions data pull)
# For numeric-looking Price (column D) values we use a leading apostrophe
# so Excel keeps them as literal text (matches original inlineStr cells)
# instead of auto-converting to a number (which would drop things like
# trailing zeros or turn "26.487.55" style multi-dot prices into errors).

$ws.Range('D2').Value = "'26.487.55"
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = "'1.731.91"
$ws.Range('E3').Value = '  -0.66%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = "'247.43"
$ws.Range('E5').Value = '  +0.42%  '

$ws.Range('D7').Value = "'0.4877"
$ws.Range('E7').Value = '  +1.04%  '

$ws.Range('D8').Value = "'0.2671"
$ws.Range('E8').Value = '  -0.83%  '

$ws.Range('D9').Value = "'0.06227"
$ws.Range('E9').Value = '  -0.55%  '

$ws.Range('D10').Value = "'1.735.77"
$ws.Range('E10').Value = '  -0.45%  '

$ws.Range('D11').Value = "'0.07066"
$ws.Range('E11').Value = '  -0.99%  '

$ws.Range('D12').Value = "'15.65"
$ws.Range('E12').Value = '  -1.62%  '

$ws.Range('D13').Value = "'4.661"
$ws.Range('E13').Value = '  +3.03%  '

$ws.Range('D14').Value = "'0.6098"
$ws.Range('E14').Value = '  -2.46%  '

$ws.Range('D15').Value = "'77.46"
$ws.Range('E15').Value = '  -0.04%  '

$ws.Range('D17').Value = "'26.492.43"
$ws.Range('E17').Value = '  -0.42%  '

$ws.Range('D18').Value = "'1.0000"
$ws.Range('E18').Value = '  +0.02%  '

$ws.Range('D19').Value = "'0.000007170"
$ws.Range('E19').Value = '  +3.65%  '

$ws.Range('D20').Value = "'11.51"
$ws.Range('E20').Value = '  -2.31%  '

$ws.Range('D21').Value = "'1.959.62"

$ws.Range('D22').Value = "'4.528"
$ws.Range('E22').Value = '  -2.11%  '

$ws.Range('D23').Value = "'8.781"
$ws.Range('E23').Value = '  -1.20%  '

$ws.Range('D24').Value = "'5.261"
$ws.Range('E24').Value = '  -2.08%  '

$ws.Range('D25').Value = "'139.29"
$ws.Range('E25').Value = '  +2.38%  '

$ws.Range('D26').Value = "'15.43"
$ws.Range('E26').Value = '  +0.33%  '

$ws.Range('D27').Value = "'1.776"
$ws.Range('E27').Value = '  -2.28%  '

$ws.Range('D28').Value = "'108.11"
$ws.Range('E28').Value = '  +1.24%  '

$ws.Range('D29').Value = "'1.402"
$ws.Range('E29').Value = '  -2.13%  '

$ws.Range('D30').Value = "'3.974"

$ws.Range('D31').Value = "'0.08036"
$ws.Range('E31').Value = '  +1.84%  '

$ws.Range('D32').Value = "'3.700"
$ws.Range('E32').Value = '  -1.18%  '

$ws.Range('E33').Value = '  -0.25%  '

$ws.Range('D34').Value = "'0.9999"

$ws.Range('D35').Value = "'2.615"
$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('E36').Value = '  +0.43%  '

$ws.Range('D37').Value = "'0.6369"
$ws.Range('E37').Value = '  -0.73%  '

$ws.Range('D38').Value = "'0.8966"
$ws.Range('E38').Value = '  -4.05%  '

$ws.Range('D39').Value = "'2.019"
$ws.Range('E39').Value = '  +1.49%  '

$ws.Range('D40').Value = "'2.398"
$ws.Range('E40').Value = '  -1.52%  '

$ws.Range('D41').Value = "'1.002"
$ws.Range('E41').Value = '  -0.19%  '

$ws.Range('E42').Value = '  -0.24%  '

$ws.Range('D43').Value = "'101.49"
$ws.Range('E43').Value = '  -11.06%  '

$ws.Range('D44').Value = "'5.446"
$ws.Range('E44').Value = '  -6.01%  '

$ws.Range('D45').Value = "'0.3891"
$ws.Range('E45').Value = '  -0.80%  '

$ws.Range('D46').Value = "'6.967"
$ws.Range('E46').Value = '  +3.11%  '

$ws.Range('D47').Value = "'0.1184"
$ws.Range('E47').Value = '  -2.86%  '

$ws.Range('D48').Value = "'0.05385"
$ws.Range('E48').Value = '  +0.89%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = "'30.63"
$ws.Range('E49').Value = '  -0.53%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'7.825"
$ws.Range('E50').Value = '  -1.47%  '

$ws.Range('E51').Value = '  -1.16%  '
